# Populate "Đơn sale chính" sheet (sheet 1) with the commission report data
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Header row
$ws1.Range("A1").Value = "Tiền tố"
$ws1.Range("B1").Value = "Mã dịch vụ"
$ws1.Range("C1").Value = "Ngày thực hiện"
$ws1.Range("D1").Value = "Cơ sở"
$ws1.Range("E1").Value = "Khách hàng"
$ws1.Range("F1").Value = "Nguồn khách"
$ws1.Range("G1").Value = "Tên dịch vụ"
$ws1.Range("H1").Value = "Đơn giá gốc"
$ws1.Range("I1").Value = "Sale phụ"
$ws1.Range("J1").Value = "Upsale"
$ws1.Range("K1").Value = "Đơn giá"
$ws1.Range("L1").Value = "Đã thanh toán"
$ws1.Range("M1").Value = "Tỉ lệ chiết khấu sale chính"
$ws1.Range("N1").Value = "Chiết khấu sale chính"

# Row 2
$ws1.Range("A2").Value = "HD-LUXURY"
$ws1.Range("B2").Value = 616
$ws1.Range("C2").Value = "'08-02-2024"
$ws1.Range("D2").Value = "LONG XUYÊN"
$ws1.Range("E2").Value = "Chị duyên"
$ws1.Range("F2").Value = "Khách cũ giới thiệu"
$ws1.Range("G2").Value = "Cắt mí"
$ws1.Range("H2").Value = 8000000
$ws1.Range("I2").Value = $null
$ws1.Range("J2").Value = $null
$ws1.Range("K2").Value = 8000000
$ws1.Range("L2").Value = 8000000
$ws1.Range("M2").Value = 0.1
$ws1.Range("N2").Value = 800000

# Row 3
$ws1.Range("A3").Value = "HD-LUXURY"
$ws1.Range("B3").Value = 617
$ws1.Range("C3").Value = "'08-02-2024"
$ws1.Range("D3").Value = "LONG XUYÊN"
$ws1.Range("E3").Value = "Cô tú"
$ws1.Range("F3").Value = "Khách cũ"
$ws1.Range("G3").Value = "Nâng cung chân mày"
$ws1.Range("H3").Value = 4000000
$ws1.Range("I3").Value = $null
$ws1.Range("J3").Value = $null
$ws1.Range("K3").Value = 4000000
$ws1.Range("L3").Value = 3000000
$ws1.Range("M3").Value = 0.1
$ws1.Range("N3").Value = 300000

# Row 4 - totals
$ws1.Range("A4").Value = "Tổng"
$ws1.Range("B4").Value = 2
$ws1.Range("C4").Value = ""
$ws1.Range("D4").Value = ""
$ws1.Range("E4").Value = ""
$ws1.Range("F4").Value = ""
$ws1.Range("G4").Value = ""
$ws1.Range("H4").Value = 12000000
$ws1.Range("I4").Value = ""
$ws1.Range("J4").Value = 0
$ws1.Range("K4").Value = 12000000
$ws1.Range("L4").Value = 11000000
$ws1.Range("M4").Value = 0
$ws1.Range("N4").Value = 1100000

# Update "Lương" sheet (sheet 2) payroll figures for LONG XUYÊN / totals
$ws2.Range("B12").Value = 2
$ws2.Range("B13").Value = 70000
$ws2.Range("B14").Value = 571428.5714285715
$ws2.Range("B15").Value = 1100000
$ws2.Range("B33").Value = 1741428.571428571
$ws2.Range("A35").Value = "Tổng lương tại HỆ THỐNG"
$ws2.Range("B35").Value = 1741428.571428571
